# Updated cryptos list on Fri Aug  2 04:49:29 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row against freshly scraped data. Rows 44/45 additionally swap
# ranking positions between VeChain and InjectiveProtocol.
#
# Column D values are assigned with a leading apostrophe so Excel keeps
# numeric-looking strings (e.g. "573.01") stored as text -- matching the
# worksheet's existing string convention -- instead of silently coercing
# them to the Number type; the Style is then reset to "Normal" so the
# apostrophe's "quote prefix" formatting marker isn't left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.517.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "'3.166.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'573.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").Value = "'164.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -3.43%  "
$ws.Range("E9").Value = "  -1.78%  "
$ws.Range("D10").Value = "'6.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("D11").Value = "'0.385"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "'3.722.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").Value = "'64.534.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "'3.166.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "'406.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").Value = "'12.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("D24").Value = "'0.486"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  -3.04%  "
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("E27").Value = "  +2.37%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'21.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").Value = "'6.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("D34").Value = "'156.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").Value = "'2.689.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").Value = "'24.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("D39").Value = "'4.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").Value = "'0.697"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'5.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("D43").Value = "'291.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.14%  "

# Rows 44/45: VeChain and InjectiveProtocol swapped ranking positions.
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0258"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'21.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.72%  "

$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "'0.0985"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("E48").Value = "  -5.22%  "
$ws.Range("D49").Value = "'10.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").Value = "'5.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").Value = "'0.879"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.31%  "
